$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "hussein2492024055414",
    "hussein2492024055454",
    "hussein2492024055715",
    "hussein2492024055813",
    "hussein2492024060011",
    "hussein2492024060102",
    "hussein2492024060845",
    "hussein2492024060929"
)

$row = 9
foreach ($val in $values) {
    $ws.Cells.Item($row, 1).Value = $val
    $ws.Cells.Item($row, 2).Value = "Test@123"
    $row = $row + 1
}
